$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Row" (column E) and "Column" (column F) values by 1 for all data rows (2-25)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 5).Value2 + 1
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 6).Value2 + 1
}

# Auto-fit column B width (Best Fit), matching the new <cols> entry in the sheet
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Select entire row 1 (matches the new selection: activeCell J8, sqref A1:XFD1048576)
$ws.Rows.Item(1).Select() | Out-Null
